$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39-43 down to 40-44.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new data record.
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 44783
$ws.Range("D39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 100112040
$ws.Range("G39").Value = "Cilantro"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 150
$ws.Range("K39").Value = 15000
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = 15000
$ws.Range("N39").Value = '$/caja 36 atados'
$ws.Range("O39").Value = "Provincia de Quillota"
$ws.Range("P39").Value = 417
$ws.Range("Q39").Value = 36
$ws.Range("R39").Value = "Hortaliza"
